$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

# Sheet "展览" (exhibition) - update "想去人数" (want-to-go count) column F
$ws1.Range("F2").Value = 130
$ws1.Range("F3").Value = 1305
$ws1.Range("F4").Value = 1110
$ws1.Range("F5").Value = 994
$ws1.Range("F6").Value = 1769
$ws1.Range("F7").Value = 548
$ws1.Range("F16").Value = 153
$ws1.Range("F21").Value = 130
$ws1.Range("F22").Value = 658
$ws1.Range("F23").Value = 26
$ws1.Range("F27").Value = 866
$ws1.Range("F29").Value = 149
$ws1.Range("F31").Value = 263

# Sheet "演出" (performance) - update column F
$ws2.Range("F5").Value = 13
$ws2.Range("F7").Value = 249
$ws2.Range("F10").Value = 619
$ws2.Range("F11").Value = 118

# Sheet "全部类型" (all types) - update column F
$ws4.Range("F3").Value = 130
$ws4.Range("F4").Value = 1305
$ws4.Range("F5").Value = 1110
$ws4.Range("F6").Value = 994
$ws4.Range("F7").Value = 1769
$ws4.Range("F8").Value = 548
$ws4.Range("F18").Value = 153
$ws4.Range("F24").Value = 13
$ws4.Range("F27").Value = 249
$ws4.Range("F28").Value = 249
$ws4.Range("F29").Value = 130
$ws4.Range("F30").Value = 658
$ws4.Range("F31").Value = 26
$ws4.Range("F33").Value = 146
$ws4.Range("F35").Value = 866
$ws4.Range("F39").Value = 149
$ws4.Range("F41").Value = 263
$ws4.Range("F42").Value = 619
$ws4.Range("F43").Value = 118
